$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (the empty paragraph just
#    before the final paragraph, near "OUTPUT"). It needs to move to a
#    new location below, and bookmark names must stay unique, so the
#    old one is deleted first.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) Merge the two runs that make up
#       "Check to see if the " + "value is less than or equal to 31 and more than or equal to 28"
#    into a single run with the combined text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Check to see if the value is less than or equal to 31 and more than or equal to 28",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Check to see if the value is less than or equal to 31 and more than or equal to 28",
    2)

# ------------------------------------------------------------------
# 3) Merge the five runs that make up
#       "Check to see if the " + "value is less than the # of days" +
#       " and more than " + "or equal to " + "1"
#    into a single run with the combined text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Check to see if the value is less than the # of days and more than or equal to 1",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Check to see if the value is less than the # of days and more than or equal to 1",
    2)

# ------------------------------------------------------------------
# 4) Add a new (collapsed) "_GoBack" bookmark right after the word
#    "Calendar" in the title paragraph, still inside that paragraph
#    (before the paragraph mark).
#
#    A bookmark range collapsed exactly at the end of a paragraph's
#    text tends to get normalised to the paragraph boundary, so a
#    one-character placeholder is inserted first to give the bookmark
#    an real, non-degenerate anchor right after "Calendar"; the
#    placeholder is then removed, leaving a zero-length bookmark in
#    exactly the right spot.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Calendar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter("X")
$markRange = $d.Range($r.Start, $r.Start + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$markRange.Text = ""
